$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E (soil_recovery_fert_N) and G (crop_recovery_fert_N) raw values per reviewer revisions.
# Dependent formula cells (F, H, I, J) recalculate automatically.
$ws.Range("E2").Value2 = 4.7911912186879002
$ws.Range("G2").Value2 = 10.610846932985739
$ws.Range("E3").Value2 = 4.7647932265003048
$ws.Range("G3").Value2 = 9.1940368270319652
$ws.Range("E4").Value2 = 3.8663687139083436
$ws.Range("G4").Value2 = 9.4789698342020579
$ws.Range("E5").Value2 = 6.0163689530113142
$ws.Range("G5").Value2 = 11.956347652120996
$ws.Range("E6").Value2 = 6.3018275512694348
$ws.Range("G6").Value2 = 11.768408370944474
$ws.Range("E7").Value2 = 2.892711990304726
$ws.Range("G7").Value2 = 7.9045601354051254
$ws.Range("E8").Value2 = 30.08956784103556
$ws.Range("G8").Value2 = 60.235246718236311
$ws.Range("E9").Value2 = 22.62803046654366
$ws.Range("G9").Value2 = 44.130107817365911
$ws.Range("E10").Value2 = 23.268411182746796
$ws.Range("G10").Value2 = 51.208162991584281
$ws.Range("E11").Value2 = 32.768436356065109
$ws.Range("G11").Value2 = 58.1850250485586
$ws.Range("E12").Value2 = 37.63873485135862
$ws.Range("G12").Value2 = 52.543907882010018
$ws.Range("E13").Value2 = 38.024546936253508
$ws.Range("G13").Value2 = 49.10605704479763
$ws.Range("E14").Value2 = 3.9630112141671456
$ws.Range("G14").Value2 = 7.7560169970800921
$ws.Range("E15").Value2 = 6.5062663298507912
$ws.Range("G15").Value2 = 10.070297454021794
$ws.Range("E16").Value2 = 6.6756696919221072
$ws.Range("G16").Value2 = 5.6443411320498686
$ws.Range("E17").Value2 = 6.5966546529759169
$ws.Range("G17").Value2 = 11.689580877675258
$ws.Range("E18").Value2 = 7.7353377222317059
$ws.Range("G18").Value2 = 12.48509853550938
$ws.Range("E19").Value2 = 5.9810781969894506
$ws.Range("G19").Value2 = 14.46204299171999
$ws.Range("E20").Value2 = 44.748012230149762
$ws.Range("G20").Value2 = 40.594439996860473
$ws.Range("E21").Value2 = 38.196418211562225
$ws.Range("G21").Value2 = 34.833175582539489
$ws.Range("E22").Value2 = 32.373323177341568
$ws.Range("G22").Value2 = 39.510795001228132
$ws.Range("E23").Value2 = 23.751043849513408
$ws.Range("G23").Value2 = 42.550634228085869
$ws.Range("E24").Value2 = 43.149062369212878
$ws.Range("G24").Value2 = 46.305885444057992
$ws.Range("E25").Value2 = 45.692881309553371
$ws.Range("G25").Value2 = 37.654484321044777

# E column now uses a 2-decimal numeric format (distinct from F/H/I/J which keep "0.0").
$ws.Range("E2:E25").NumberFormat = "0.00"

# New helper columns N and O added alongside M; all three get a 5-decimal format.
$ws.Range("M2:O25").NumberFormat = "0.00000"

# Match the reviewer selection left behind in the sheet view.
$ws.Range("M2:O26").Select()
